# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap Rusia/Portugal rows: Rusia moves up to row 18 (with new data),
#     Portugal moves down to row 19 (keeping its previous data).
$ws.Range("A18").Value = "Rusia"
$ws.Range("B18").Value = 18328
$ws.Range("C18").Value = 2558
$ws.Range("D18").Value = 1470
$ws.Range("E18").Value = 16710
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = 148

$ws.Range("A19").Value = "Portugal"
$ws.Range("B19").Value = 16585
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 277
$ws.Range("E19").Value = 15804
$ws.Range("F19").Value = 228
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 504

# --- Update Australia row (row 32)
$ws.Range("B32").Value = 6359
$ws.Range("C32").Value = 46
$ws.Range("D32").Value = 3494
$ws.Range("E32").Value = 2804
$ws.Range("F32").Value = 79

# --- Update Moldavia row (row 60)
$ws.Range("D60").Value = 107
$ws.Range("E60").Value = 1522

# --- Update Estonia row (row 66)
$ws.Range("B66").Value = 1332
$ws.Range("C66").Value = 23
$ws.Range("D66").Value = 102
$ws.Range("E66").Value = 1203
$ws.Range("F66").Value = 9
$ws.Range("G66").Value = 2
$ws.Range("H66").Value = 27

# --- Update Armenia row (row 72)
$ws.Range("B72").Value = 1039
$ws.Range("C72").Value = 26
$ws.Range("D72").Value = 211
$ws.Range("E72").Value = 814
$ws.Range("G72").Value = 1
$ws.Range("H72").Value = 14

# --- Update Bosnia y Herzegovina row (row 73)
$ws.Range("B73").Value = 1018
$ws.Range("C73").Value = 9
$ws.Range("D73").Value = 206
$ws.Range("E73").Value = 773

# --- Update the "Datos actualizados" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 09:52"
